# Unit 1 Computer Science 9-1 - Complete Presentation Slides
# Adds the "Decomposition and Abstraction" unit content:
#  - fills in the previously-blank "4. Decomposition and Abstraction" section
#    divider (existing slide 20) and its follow-up definition slide (slide 21)
#  - appends six new slides (22-27) covering levels of abstraction and the
#    noughts-and-crosses worked example

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 20 : "4. Decomposition and Abstraction" section title
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$s20.Shapes.Item(1).TextFrame.TextRange.Text = "4. Decomposition and Abstraction"
$s20.Shapes.Item(2).TextFrame.TextRange.Text = "Analyze a program, investigate requirements (inputs, outputs, processing, initialization) and design solutions.`rDecompose a problem into smaller sub-problems`rUnderstand how abstraction can be used effectively to model aspects of the real world`rProgram abstractions of real-world example."

# ---------------------------------------------------------------------------
# Slide 21 : "Definition of decomposition and abstraction"
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$s21.Shapes.Item(1).TextFrame.TextRange.Text = "Definition of decomposition and abstraction"
$s21.Shapes.Item(1).TextFrame.AutoSize = 2
$s21.Shapes.Item(2).TextFrame.TextRange.Text = "Decomposition is when breaking a problem down inro smaller, more manageable parts, which are then easier to solve. `rAbstraction is when the process of removing or hiding necessary detail so that only the important points remain`r"

# ---------------------------------------------------------------------------
# Slide 22 (new) : Design and create noughts and crosses game
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Add(22, 2)
$s22.Shapes.Item(1).TextFrame.TextRange.Text = "Design and create noughts and crosses game"
$s22.Shapes.Item(1).TextFrame.AutoSize = 2
$s22.Shapes.Item(2).TextFrame.TextRange.Text = "The goal in designing algorithms is the fact that we shall solve sub-programs. In order to do that we first need to clarify the sub-programs. `rIn this game, the first step is to design an interface showing the 3 x 3 grid. `rThe second step, is to keep track of which squares have been selected by X and 0 and which are free. `rThe third step is how the computer will decide which square to select. `rThe last one is how the computer will decide when the game is over and who has won. "
$s22.Shapes.Item(2).TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 23 (new) : Levels of abstraction
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Add(23, 2)
$s23.Shapes.Item(1).TextFrame.TextRange.Text = "Levels of abstraction"
$s23.Shapes.Item(2).TextFrame.TextRange.Text = "There are different levels or types of abstraction. The higher the level of abstraction, the less detail is required. We use abstraction all the time in accomplishing everyday tasks. `rWhen programmers write the print command, they do not have to bother about all of the details of how this will be accomplished. They are removed from them. They are at a certain level of abstraction. `rA driver turning the ignition key to start a car does not have to understand how the engine works or  how the spark to ignite the petrol is generated. It just happens and they can simply drive the car. That is abstraction in less detail. "
$s23.Shapes.Item(2).TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 24 (new) : Noughts and crosses game abstraction
# ---------------------------------------------------------------------------
$s24 = $p.Slides.Add(24, 2)
$s24.Shapes.Item(1).TextFrame.TextRange.Text = "Noughts and crosses game abstraction"
$s24.Shapes.Item(2).TextFrame.TextRange.Text = "The computer goes first. Then the user. This continues until either one wins, or all of the squares have been used. It shall generated through inputs and outputs. `rStart the game. `rEntries for the user. `rSelect a new game or finish. `rA message to inform the user when it is their turn. `rA message to inform the user if they try to select a square that has already been used. "

# ---------------------------------------------------------------------------
# Slide 25 (new) : ...and more on noughts and crosses Part I
# ---------------------------------------------------------------------------
$s25 = $p.Slides.Add(25, 2)
$s25.Shapes.Item(1).TextFrame.TextRange.Text = [char]0x2026 + "and more on noughts and crosses Part I"
$s25.Shapes.Item(2).TextFrame.TextRange.Text = "A message to inform the user if the game is a draw. `rA message to inform the user if they or the computer has won. `rA message to ask the user if they want to play another game or want to finish `rTherefore, we move on to processing and initialization`rSet up the grid with nine squares. `rInitialise all variables to a start value. `rDecide which square the computer will select. "

# ---------------------------------------------------------------------------
# Slide 26 (new) : ...and more noughts and crosses Part II
# ---------------------------------------------------------------------------
$s26 = $p.Slides.Add(26, 2)
$s26.Shapes.Item(1).TextFrame.TextRange.Text = [char]0x2026 + "and more noughts and crosses Part II"
$s26.Shapes.Item(2).TextFrame.TextRange.Text = "Allow the user to select a square. `rCheck if the user has selected an already used square. `rCheck if the computer or the user has won. `rCheck if all squares have been used and the game is a draw. `rAllow the user to select a new game or finish. `rHomework : Exercise : Code the aforementioned program in Python or Java. "

# ---------------------------------------------------------------------------
# Slide 27 (new) : blank placeholder slide
# ---------------------------------------------------------------------------
$s27 = $p.Slides.Add(27, 2)
